# Update cryptos list - generated from diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / already-non-numeric-looking cell updates ---
$ws.Range("D2").Value = "67.481.70"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.628.91"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("D9").Value = "2.628.04"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").Value = "3.106.90"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "67.141.07"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "2.630.92"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("E20").Value = "  +4.44%  "
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("E23").Value = "  -2.97%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  -4.82%  "
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").Value = "2.761.54"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -3.72%  "
$ws.Range("D46").Value = "0.0₆0297"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E51").Value = "  -1.38%  "

# --- Numeric-looking price strings: force text storage so they keep their exact
#     string form (e.g. "5.22") instead of being parsed as numbers ---
$forceTextRanges = @("D5","D6","D8","D12","D14","D16","D19","D20","D21","D22","D23","D25","D26","D27","D30","D31","D32","D33","D34","D38","D39","D40","D41","D42","D43","D45","D47","D48","D49","D50","D51")
foreach ($r in $forceTextRanges) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D5").Value = "594.17"
$ws.Range("D6").Value = "168.59"
$ws.Range("D8").Value = "0.534"
$ws.Range("D12").Value = "0.365"
$ws.Range("D14").Value = "27.68"
$ws.Range("D16").Value = "0.0000183"
$ws.Range("D19").Value = "12.04"
$ws.Range("D20").Value = "8.02"
$ws.Range("D21").Value = "356.77"
$ws.Range("D22").Value = "4.32"
$ws.Range("D23").Value = "4.67"
$ws.Range("D25").Value = "1.93"
$ws.Range("D26").Value = "10.32"
$ws.Range("D27").Value = "69.61"
$ws.Range("D30").Value = "0.0000101"
$ws.Range("D31").Value = "545.13"
$ws.Range("D32").Value = "7.92"
$ws.Range("D33").Value = "1.35"
$ws.Range("D34").Value = "1.90"
$ws.Range("D38").Value = "156.37"
$ws.Range("D39").Value = "19.03"
$ws.Range("D40").Value = "0.366"
$ws.Range("D41").Value = "1.82"
$ws.Range("D42").Value = "5.22"
$ws.Range("D43").Value = "18.20"
$ws.Range("D45").Value = "2.43"
$ws.Range("D47").Value = "152.87"
$ws.Range("D48").Value = "0.580"
$ws.Range("D49").Value = "3.79"
$ws.Range("D50").Value = "1.70"
$ws.Range("D51").Value = "0.0770"

# Reset number format back to General/Normal so no stray text-format style lingers
foreach ($r in $forceTextRanges) {
    $ws.Range($r).Style = "Normal"
}
